$wb = $excel.ActiveWorkbook

# ALC row 98
$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H98").Value = 1914.75
$ws_ALC.Range("I98").Value = 1757.7693
$ws_ALC.Range("J98").Value = 2595
$ws_ALC.Range("K98").Value = 1757.7693
$ws_ALC.Range("L98").Value = 2595
$ws_ALC.Range("M98").Value = -259.7692999999999
$ws_ALC.Range("N98").Value = -5591

# ALC row 112
$ws_ALC.Range("H112").Value = 2700.2727
$ws_ALC.Range("J112").Value = 3872.0715
$ws_ALC.Range("L112").Value = 11616.2145
$ws_ALC.Range("N112").Value = -13832.2145

# ALC row 113
$ws_ALC.Range("H113").Value = 2503.68
$ws_ALC.Range("I113").Value = 3000
$ws_ALC.Range("J113").Value = 2483
$ws_ALC.Range("K113").Value = 3000
$ws_ALC.Range("L113").Value = 2483
$ws_ALC.Range("M113").Value = 254
$ws_ALC.Range("N113").Value = -8991

# ALC row 122
$ws_ALC.Range("H122").Value = 1914.75
$ws_ALC.Range("I122").Value = 1757.7693
$ws_ALC.Range("J122").Value = 2595
$ws_ALC.Range("K122").Value = 5273.3079
$ws_ALC.Range("L122").Value = 7785
$ws_ALC.Range("M122").Value = -2823.3079
$ws_ALC.Range("N122").Value = -12685

# ALC row 129
$ws_ALC.Range("H129").Value = 878.2432
$ws_ALC.Range("I129").Value = 498.5
$ws_ALC.Range("J129").Value = 899.9429
$ws_ALC.Range("K129").Value = 1495.5
$ws_ALC.Range("L129").Value = 2699.8287
$ws_ALC.Range("M129").Value = 3504.5
$ws_ALC.Range("N129").Value = -12699.8287

# ALC row 138
$ws_ALC.Range("H138").Value = 822152.5
$ws_ALC.Range("J138").Value = 1280974.4
$ws_ALC.Range("L138").Value = 3842923.2
$ws_ALC.Range("N138").Value = -3853203.2

# ARM row 5
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H5").Value = 169
$ws_ARM.Range("I5").Value = 169
$ws_ARM.Range("K5").Value = 169
$ws_ARM.Range("M5").Value = -57

# ARM row 32
$ws_ARM.Range("H32").Value = 3327.762
$ws_ARM.Range("I32").Value = 3459.025
$ws_ARM.Range("K32").Value = 3459.025
$ws_ARM.Range("M32").Value = -3172.025

# ARM row 45
$ws_ARM.Range("H45").Value = 1908.7858
$ws_ARM.Range("J45").Value = 1831.6666
$ws_ARM.Range("L45").Value = 1831.6666
$ws_ARM.Range("N45").Value = -2585.6666

# ARM row 61
$ws_ARM.Range("H61").Value = 1534.5834
$ws_ARM.Range("I61").Value = 1055.6666
$ws_ARM.Range("K61").Value = 1055.6666
$ws_ARM.Range("M61").Value = -843.6666

# ARM row 122
$ws_ARM.Range("H122").Value = 1489.75
$ws_ARM.Range("I122").Value = 1220.8
$ws_ARM.Range("J122").Value = 1938
$ws_ARM.Range("K122").Value = 3662.4
$ws_ARM.Range("L122").Value = 5814
$ws_ARM.Range("M122").Value = -1212.4
$ws_ARM.Range("N122").Value = -10714

# ARM row 128
$ws_ARM.Range("H128").Value = 69740
$ws_ARM.Range("J128").Value = 69740
$ws_ARM.Range("L128").Value = 69740
$ws_ARM.Range("N128").Value = -79700

# ARM row 132
$ws_ARM.Range("H132").Value = 2752.7778
$ws_ARM.Range("I132").Value = 2491.762
$ws_ARM.Range("J132").Value = 3666.3333
$ws_ARM.Range("K132").Value = 7475.286
$ws_ARM.Range("L132").Value = 10998.9999
$ws_ARM.Range("M132").Value = -4945.286
$ws_ARM.Range("N132").Value = -16058.9999

# ARM row 136
$ws_ARM.Range("H136").Value = 1534.5834
$ws_ARM.Range("I136").Value = 1055.6666
$ws_ARM.Range("K136").Value = 3166.9998
$ws_ARM.Range("M136").Value = -616.9998000000001

# BSM row 4
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H4").Value = 169
$ws_BSM.Range("I4").Value = 169
$ws_BSM.Range("K4").Value = 169
$ws_BSM.Range("M4").Value = -54

# BSM row 105
$ws_BSM.Range("H105").Value = 63120610
$ws_BSM.Range("I105").Value = 100991170
$ws_BSM.Range("K105").Value = 100991170
$ws_BSM.Range("M105").Value = -100989423

# BSM row 107
$ws_BSM.Range("H107").Value = 1248.9048
$ws_BSM.Range("I107").Value = 1015.2857
$ws_BSM.Range("J107").Value = 1716.1428
$ws_BSM.Range("K107").Value = 1015.2857
$ws_BSM.Range("L107").Value = 1716.1428
$ws_BSM.Range("M107").Value = 904.7143
$ws_BSM.Range("N107").Value = -5556.1428

# BSM row 134
$ws_BSM.Range("H134").Value = 14766.125
$ws_BSM.Range("I134").Value = 2426
$ws_BSM.Range("K134").Value = 7278
$ws_BSM.Range("M134").Value = -4743

# CRP row 31
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 877.7857
$ws_CRP.Range("I31").Value = 739.1905
$ws_CRP.Range("J31").Value = 1293.5714
$ws_CRP.Range("K31").Value = 739.1905
$ws_CRP.Range("L31").Value = 1293.5714
$ws_CRP.Range("M31").Value = -444.1905
$ws_CRP.Range("N31").Value = -1883.5714

# CRP row 34
$ws_CRP.Range("H34").Value = 877.7857
$ws_CRP.Range("I34").Value = 739.1905
$ws_CRP.Range("J34").Value = 1293.5714
$ws_CRP.Range("K34").Value = 739.1905
$ws_CRP.Range("L34").Value = 1293.5714
$ws_CRP.Range("M34").Value = -537.1905
$ws_CRP.Range("N34").Value = -1697.5714

# CUL row 68
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H68").Value = 1407.3125
$ws_CUL.Range("I68").Value = 707.9167
$ws_CUL.Range("J68").Value = 1826.95
$ws_CUL.Range("K68").Value = 2123.7501
$ws_CUL.Range("L68").Value = 5480.85
$ws_CUL.Range("M68").Value = -1312.7501
$ws_CUL.Range("N68").Value = -7102.85

# CUL row 71
$ws_CUL.Range("H71").Value = 1407.3125
$ws_CUL.Range("I71").Value = 707.9167
$ws_CUL.Range("J71").Value = 1826.95
$ws_CUL.Range("K71").Value = 6371.2503
$ws_CUL.Range("L71").Value = 16442.55
$ws_CUL.Range("M71").Value = -2315.2503
$ws_CUL.Range("N71").Value = -24554.55

# CUL row 81
$ws_CUL.Range("H81").Value = 2570.3076
$ws_CUL.Range("I81").Value = 1196.6
$ws_CUL.Range("J81").Value = 2897.3809
$ws_CUL.Range("K81").Value = 3589.8
$ws_CUL.Range("L81").Value = 8692.1427
$ws_CUL.Range("M81").Value = -2466.8
$ws_CUL.Range("N81").Value = -10938.1427

# CUL row 84
$ws_CUL.Range("H84").Value = 2570.3076
$ws_CUL.Range("I84").Value = 1196.6
$ws_CUL.Range("J84").Value = 2897.3809
$ws_CUL.Range("K84").Value = 10769.4
$ws_CUL.Range("L84").Value = 26076.4281
$ws_CUL.Range("M84").Value = -5153.4
$ws_CUL.Range("N84").Value = -37308.4281

# CUL row 107
$ws_CUL.Range("H107").Value = 5082
$ws_CUL.Range("J107").Value = 9949
$ws_CUL.Range("L107").Value = 29847
$ws_CUL.Range("N107").Value = -33687

# CUL row 113
$ws_CUL.Range("H113").Value = 605.1667
$ws_CUL.Range("I113").Value = 600
$ws_CUL.Range("J113").Value = 605.4706
$ws_CUL.Range("K113").Value = 1800
$ws_CUL.Range("L113").Value = 1816.4118
$ws_CUL.Range("M113").Value = 370
$ws_CUL.Range("N113").Value = -6156.4118

# CUL row 131
$ws_CUL.Range("H131").Value = 18519834
$ws_CUL.Range("J131").Value = 1383.625
$ws_CUL.Range("L131").Value = 4150.875
$ws_CUL.Range("N131").Value = -14230.875

# CUL row 132
$ws_CUL.Range("H132").Value = 0
$ws_CUL.Range("I132").Value = 0
$ws_CUL.Range("J132").Value = 0
$ws_CUL.Range("K132").Value = 0
$ws_CUL.Range("L132").Value = 0
$ws_CUL.Range("M132").ClearContents()
$ws_CUL.Range("N132").ClearContents()

# GSM row 113
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H113").Value = 2078.3333
$ws_GSM.Range("I113").Value = 1206.6666
$ws_GSM.Range("K113").Value = 1206.6666
$ws_GSM.Range("M113").Value = 963.3334

# GSM row 126
$ws_GSM.Range("H126").Value = 2040.6111
$ws_GSM.Range("I126").Value = 1629.7273
$ws_GSM.Range("K126").Value = 4889.1819
$ws_GSM.Range("M126").Value = -2419.1819

# GSM row 132
$ws_GSM.Range("H132").Value = 2536.037
$ws_GSM.Range("I132").Value = 2128.0588
$ws_GSM.Range("K132").Value = 6384.176399999999
$ws_GSM.Range("M132").Value = -3854.176399999999

# LTW row 40
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H40").Value = 2238.7144
$ws_LTW.Range("I40").Value = 1984.8182
$ws_LTW.Range("J40").Value = 3169.6667
$ws_LTW.Range("K40").Value = 1984.8182
$ws_LTW.Range("L40").Value = 3169.6667
$ws_LTW.Range("M40").Value = -1848.8182
$ws_LTW.Range("N40").Value = -3441.6667

# LTW row 92
$ws_LTW.Range("H92").Value = 17750
$ws_LTW.Range("J92").Value = 17750
$ws_LTW.Range("L92").Value = 17750
$ws_LTW.Range("N92").Value = -22742

# WVR row 96
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H96").Value = 1948.3529
$ws_WVR.Range("I96").Value = 2003.5
$ws_WVR.Range("J96").Value = 1816
$ws_WVR.Range("K96").Value = 2003.5
$ws_WVR.Range("L96").Value = 1816
$ws_WVR.Range("M96").Value = -630.5
$ws_WVR.Range("N96").Value = -4562

# WVR row 136
$ws_WVR.Range("H136").Value = 1989
$ws_WVR.Range("I136").Value = 1779.4
$ws_WVR.Range("J136").Value = 2251
$ws_WVR.Range("K136").Value = 5338.200000000001
$ws_WVR.Range("L136").Value = 6753
$ws_WVR.Range("M136").Value = -2788.200000000001
